$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new header cells H1:J1 ("pv", "cv", "fv"), copying the
#     existing header formatting (style) from G1 so they match the other
#     header cells in row 1. ---
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "pv"
$ws.Range("I1").Value = "cv"
$ws.Range("J1").Value = "fv"

# --- Fill in the new data columns H (pv), I (cv), J (fv) for rows 2-24 ---
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 1

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 0

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 14
$ws.Range("J4").Value = 3

$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 19
$ws.Range("J5").Value = 0

$ws.Range("H6").Value = 14
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5

$ws.Range("H7").Value = 18
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3

$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 2

$ws.Range("H9").Value = 1.7
$ws.Range("I9").Value = 95.6
$ws.Range("J9").Value = 2.7

$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 92.8
$ws.Range("J10").Value = 3.1

$ws.Range("H11").Value = 35.5
$ws.Range("I11").Value = 2.7
$ws.Range("J11").Value = 61.9

$ws.Range("H12").Value = 21.8
$ws.Range("I12").Value = 29.1
$ws.Range("J12").Value = 49.1

$ws.Range("H13").Value = 2.3
$ws.Range("I13").Value = 47.1
$ws.Range("J13").Value = 50.6

$ws.Range("H14").Value = 12
$ws.Range("I14").Value = 73.6
$ws.Range("J14").Value = 14.4

$ws.Range("H15").Value = 14.3
$ws.Range("I15").Value = 32.1
$ws.Range("J15").Value = 53.6

$ws.Range("H16").Value = 7.4
$ws.Range("I16").Value = 51.4
$ws.Range("J16").Value = 41.3

$ws.Range("H17").Value = 11.3
$ws.Range("I17").Value = 38.1
$ws.Range("J17").Value = 50.6

$ws.Range("H18").Value = 23.2
$ws.Range("I18").Value = 31.9
$ws.Range("J18").Value = 44.8

$ws.Range("H19").Value = 15.6
$ws.Range("I19").Value = 6.8
$ws.Range("J19").Value = 77.6

$ws.Range("H20").Value = 20.4
$ws.Range("I20").Value = 65.7
$ws.Range("J20").Value = 13.9

$ws.Range("H21").Value = 14
$ws.Range("I21").Value = 31.1
$ws.Range("J21").Value = 54.9

$ws.Range("H22").Value = 8.3
$ws.Range("I22").Value = 58.3
$ws.Range("J22").Value = 33.3

$ws.Range("H23").Value = 19.2
$ws.Range("I23").Value = 44.8
$ws.Range("J23").Value = 36.1

$ws.Range("H24").Value = 24.2
$ws.Range("I24").Value = 20.9
$ws.Range("J24").Value = 54.8

# --- Update the view: active selection moves to J26 and the sheet is
#     scrolled down/right so that row 6 / column B is the top-left
#     visible cell. ---
[void]$ws.Range("J26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 2

# --- Resize the saved workbook window metadata. ---
$win.Width = 19200
$win.Height = 7050
